# Add 5 more super-adventure-game entries (6 new rows, continuing 1986 releases)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1986, "吸血鬼の洞窟",   "Crypt of the Vampire",      "Tokyo Sogensha", "vampire_cave.jpg"),
    @(1986, "シャドー砦の魔王", "The Lord of Shadow Keep",   "Tokyo Sogensha", "demon_king_of_fort_shadow.jpg"),
    @(1986, "炎の神殿 ",       "The Temple of Flame",       "Tokyo Sogensha", "temple_of_the_flame.jpg"),
    @(1986, "失われた魂の城",   "Castle of Lost Souls",      "Tokyo Sogensha", "castle_of_lost_souls.jpg"),
    @(1986, "ドラゴンの目",     "The Eye of the Dragon",     "Tokyo Sogensha", "eye_of_the_dragon.jpg"),
    @(1986, "ファラオの呪い",   "Curse of the Pharaoh",      "Tokyo Sogensha", "curse_of_the_pharaoh.jpg")
)

$startRow = 31
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("E37").Select()
